# Auto-generated cell updates applying scheduled-runner market data refresh
# to the Kujata Profits workbook (columns H-N per job-class sheet).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 471.42856
$ws.Range("I8").Value = 471.42856
$ws.Range("K8").Value = 1414.28568
$ws.Range("M8").Value = -1275.28568
# Row 18
$ws.Range("H18").Value = 1999
$ws.Range("I18").Value = 2198
$ws.Range("J18").Value = 1667.3334
$ws.Range("K18").Value = 2198
$ws.Range("L18").Value = 1667.3334
$ws.Range("M18").Value = -1914
$ws.Range("N18").Value = -2235.3334
# Row 28
$ws.Range("H28").Value = 144.26666
$ws.Range("I28").Value = 146.5
$ws.Range("J28").Value = 139.8
$ws.Range("K28").Value = 146.5
$ws.Range("L28").Value = 139.8
$ws.Range("M28").Value = 338.5
$ws.Range("N28").Value = -1109.8
# Row 62
$ws.Range("H62").Value = 9999
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 9999
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -9375
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 9999
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 49995
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -46875
$ws.Range("N65").ClearContents()
# Row 70
$ws.Range("H70").Value = 1591.3572
$ws.Range("I70").Value = 1477.3334
$ws.Range("J70").Value = 1796.6
$ws.Range("K70").Value = 4432.0002
$ws.Range("L70").Value = 5389.799999999999
$ws.Range("M70").Value = -4162.0002
$ws.Range("N70").Value = -5929.799999999999
# Row 73
$ws.Range("H73").Value = 1591.3572
$ws.Range("I73").Value = 1477.3334
$ws.Range("J73").Value = 1796.6
$ws.Range("K73").Value = 4432.0002
$ws.Range("L73").Value = 5389.799999999999
$ws.Range("M73").Value = -3496.0002
$ws.Range("N73").Value = -7261.799999999999
# Row 92
$ws.Range("H92").Value = 494.5
$ws.Range("I92").Value = 340.6875
$ws.Range("K92").Value = 340.6875
$ws.Range("M92").Value = 907.3125
# Row 98
$ws.Range("H98").Value = 3905.625
$ws.Range("I98").Value = 2489.3333
$ws.Range("J98").Value = 25150
$ws.Range("K98").Value = 2489.3333
$ws.Range("L98").Value = 25150
$ws.Range("M98").Value = -991.3332999999998
$ws.Range("N98").Value = -28146
# Row 106
$ws.Range("H106").Value = 8700.883
$ws.Range("I106").Value = 9807.786
$ws.Range("K106").Value = 9807.786
$ws.Range("M106").Value = -9176.786
# Row 107
$ws.Range("H107").Value = 1418.3846
$ws.Range("I107").Value = 1289.5294
$ws.Range("J107").Value = 1661.7778
$ws.Range("K107").Value = 1289.5294
$ws.Range("L107").Value = 1661.7778
$ws.Range("M107").Value = 630.4706000000001
$ws.Range("N107").Value = -5501.7778
# Row 111
$ws.Range("H111").Value = 1017.9091
$ws.Range("I111").Value = 1030.8334
$ws.Range("J111").Value = 1002.4
$ws.Range("K111").Value = 3092.5002
$ws.Range("L111").Value = 3007.2
$ws.Range("M111").Value = -25.50019999999995
$ws.Range("N111").Value = -9141.200000000001
# Row 112
$ws.Range("H112").Value = 2335.1177
$ws.Range("J112").Value = 2733.111
$ws.Range("L112").Value = 8199.332999999999
$ws.Range("N112").Value = -10415.333
# Row 122
$ws.Range("H122").Value = 3905.625
$ws.Range("I122").Value = 2489.3333
$ws.Range("J122").Value = 25150
$ws.Range("K122").Value = 7467.999899999999
$ws.Range("L122").Value = 75450
$ws.Range("M122").Value = -5017.999899999999
$ws.Range("N122").Value = -80350
# Row 138
$ws.Range("H138").Value = 1482.4
$ws.Range("I138").Value = 586.871
$ws.Range("J138").Value = 1884.7391
$ws.Range("K138").Value = 1760.613
$ws.Range("L138").Value = 5654.2173
$ws.Range("M138").Value = 3379.387
$ws.Range("N138").Value = -15934.2173
# Row 141
$ws.Range("H141").Value = 794
$ws.Range("I141").Value = 794
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2382
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2798
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 10000000
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999856
# Row 32
$ws.Range("H32").Value = 4032.0386
$ws.Range("I32").Value = 3959.7778
$ws.Range("J32").Value = 4496.5713
$ws.Range("K32").Value = 3959.7778
$ws.Range("L32").Value = 4496.5713
$ws.Range("M32").Value = -3672.7778
$ws.Range("N32").Value = -5070.5713
# Row 61
$ws.Range("H61").Value = 1482.3572
$ws.Range("I61").Value = 969.375
$ws.Range("K61").Value = 969.375
$ws.Range("M61").Value = -757.375
# Row 74
$ws.Range("H74").Value = 1763
$ws.Range("I74").Value = 755.1667
$ws.Range("J74").Value = 2434.889
$ws.Range("K74").Value = 755.1667
$ws.Range("L74").Value = 2434.889
$ws.Range("M74").Value = 118.8333
$ws.Range("N74").Value = -4182.889
# Row 77
$ws.Range("H77").Value = 1763
$ws.Range("I77").Value = 755.1667
$ws.Range("J77").Value = 2434.889
$ws.Range("K77").Value = 3775.8335
$ws.Range("L77").Value = 12174.445
$ws.Range("M77").Value = 592.1665000000003
$ws.Range("N77").Value = -20910.445
# Row 132
$ws.Range("H132").Value = 2410.724
$ws.Range("I132").Value = 2087.2273
$ws.Range("J132").Value = 3427.4285
$ws.Range("K132").Value = 6261.6819
$ws.Range("L132").Value = 10282.2855
$ws.Range("M132").Value = -3731.6819
$ws.Range("N132").Value = -15342.2855
# Row 136
$ws.Range("H136").Value = 1482.3572
$ws.Range("I136").Value = 969.375
$ws.Range("K136").Value = 2908.125
$ws.Range("M136").Value = -358.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5951
$ws.Range("I134").Value = 1086.8572
$ws.Range("J134").Value = 40000
$ws.Range("K134").Value = 3260.5716
$ws.Range("L134").Value = 120000
$ws.Range("M134").Value = -725.5715999999998
$ws.Range("N134").Value = -125070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1244.7534
$ws.Range("I31").Value = 1097.2222
$ws.Range("J31").Value = 2174.2
$ws.Range("K31").Value = 1097.2222
$ws.Range("L31").Value = 2174.2
$ws.Range("M31").Value = -802.2221999999999
$ws.Range("N31").Value = -2764.2
# Row 34
$ws.Range("H34").Value = 1244.7534
$ws.Range("I34").Value = 1097.2222
$ws.Range("J34").Value = 2174.2
$ws.Range("K34").Value = 1097.2222
$ws.Range("L34").Value = 2174.2
$ws.Range("M34").Value = -895.2221999999999
$ws.Range("N34").Value = -2578.2
# Row 58
$ws.Range("H58").Value = 1213.05
$ws.Range("J58").Value = 1840.2
$ws.Range("L58").Value = 1840.2
$ws.Range("N58").Value = -2246.2
# Row 86
$ws.Range("H86").Value = 3937023.5
$ws.Range("I86").Value = 8336735.5
$ws.Range("J86").Value = 26168.223
$ws.Range("K86").Value = 8336735.5
$ws.Range("L86").Value = 26168.223
$ws.Range("M86").Value = -8335612.5
$ws.Range("N86").Value = -28414.223
# Row 89
$ws.Range("H89").Value = 3937023.5
$ws.Range("I89").Value = 8336735.5
$ws.Range("J89").Value = 26168.223
$ws.Range("K89").Value = 41683677.5
$ws.Range("L89").Value = 130841.115
$ws.Range("M89").Value = -41678061.5
$ws.Range("N89").Value = -142073.115
# Row 98
$ws.Range("H98").Value = 40780
$ws.Range("J98").Value = 40780
$ws.Range("L98").Value = 40780
$ws.Range("N98").Value = -45272
# Row 105
$ws.Range("H105").Value = 822.375
$ws.Range("I105").Value = 779.8333
$ws.Range("J105").Value = 950
$ws.Range("K105").Value = 779.8333
$ws.Range("L105").Value = 950
$ws.Range("M105").Value = 967.1667
$ws.Range("N105").Value = -4444
# Row 134
$ws.Range("H134").Value = 1101.8788
$ws.Range("I134").Value = 947.06665
$ws.Range("K134").Value = 2841.19995
$ws.Range("M134").Value = -306.1999500000002
# Row 136
$ws.Range("H136").Value = 1213.05
$ws.Range("J136").Value = 1840.2
$ws.Range("L136").Value = 5520.6
$ws.Range("N136").Value = -10620.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 54
$ws.Range("H54").Value = 2667.6667
$ws.Range("J54").Value = 2667.6667
$ws.Range("L54").Value = 8003.000100000001
$ws.Range("N54").Value = -9121.000100000001
# Row 102
$ws.Range("H102").Value = 2982.25
$ws.Range("J102").Value = 2982.25
$ws.Range("L102").Value = 8946.75
$ws.Range("N102").Value = -13814.75
# Row 107
$ws.Range("H107").Value = 14839.857
$ws.Range("J107").Value = 17146.5
$ws.Range("L107").Value = 51439.5
$ws.Range("N107").Value = -55279.5
# Row 131
$ws.Range("H131").Value = 17244112
$ws.Range("J131").Value = 3017.25
$ws.Range("L131").Value = 9051.75
$ws.Range("N131").Value = -19131.75
# Row 136
$ws.Range("H136").Value = 2832.7144
$ws.Range("I136").Value = 1607.5
$ws.Range("J136").Value = 4466.3335
$ws.Range("K136").Value = 4822.5
$ws.Range("L136").Value = 13399.0005
$ws.Range("M136").Value = 277.5
$ws.Range("N136").Value = -23599.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4675
$ws.Range("I80").Value = 4774
$ws.Range("J80").Value = 4510
$ws.Range("K80").Value = 4774
$ws.Range("L80").Value = 4510
$ws.Range("M80").Value = -3776
$ws.Range("N80").Value = -6506
# Row 83
$ws.Range("H83").Value = 4675
$ws.Range("I83").Value = 4774
$ws.Range("J83").Value = 4510
$ws.Range("K83").Value = 23870
$ws.Range("L83").Value = 22550
$ws.Range("M83").Value = -18878
$ws.Range("N83").Value = -32534
# Row 132
$ws.Range("H132").Value = 1855.6428
$ws.Range("I132").Value = 1524.2106
$ws.Range("J132").Value = 2555.3333
$ws.Range("K132").Value = 4572.6318
$ws.Range("L132").Value = 7665.999899999999
$ws.Range("M132").Value = -2042.6318
$ws.Range("N132").Value = -12725.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1130.2727
$ws.Range("I16").Value = 1093.3
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1093.3
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -923.3
$ws.Range("N16").Value = -1840
# Row 100
$ws.Range("H100").Value = 2231.5
$ws.Range("I100").Value = 2112
$ws.Range("K100").Value = 2112
$ws.Range("M100").Value = -1571
# Row 132
$ws.Range("H132").Value = 21522.9
$ws.Range("I132").Value = 1001.23334
$ws.Range("J132").Value = 52305.4
$ws.Range("K132").Value = 3003.70002
$ws.Range("L132").Value = 156916.2
$ws.Range("M132").Value = -473.7000200000002
$ws.Range("N132").Value = -161976.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 8966863
$ws.Range("I122").Value = 11819549
$ws.Range("J122").Value = 1279
$ws.Range("K122").Value = 35458647
$ws.Range("L122").Value = 3837
$ws.Range("M122").Value = -35456197
$ws.Range("N122").Value = -8737
# Row 126
$ws.Range("H126").Value = 90910584
$ws.Range("I126").Value = 90910584
$ws.Range("K126").Value = 272731752
$ws.Range("M126").Value = -272729282
# Row 132
$ws.Range("H132").Value = 911.75
$ws.Range("I132").Value = 497.32352
$ws.Range("K132").Value = 1491.97056
$ws.Range("M132").Value = 1038.02944
